# Auto-generated script to apply market-data refresh values
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1426.2632
$ws.Range("I19").Value = 790.8570999999999
$ws.Range("J19").Value = 1796.9166
$ws.Range("K19").Value = 790.8570999999999
$ws.Range("L19").Value = 1796.9166
$ws.Range("M19").Value = -615.8570999999999
$ws.Range("N19").Value = -2146.9166
$ws.Range("H62").Value = 4359.125
$ws.Range("I62").Value = 3526.6667
$ws.Range("J62").Value = 4858.6
$ws.Range("K62").Value = 3526.6667
$ws.Range("L62").Value = 4858.6
$ws.Range("M62").Value = -2902.6667
$ws.Range("N62").Value = -6106.6
$ws.Range("H65").Value = 4359.125
$ws.Range("I65").Value = 3526.6667
$ws.Range("J65").Value = 4858.6
$ws.Range("K65").Value = 17633.3335
$ws.Range("L65").Value = 24293
$ws.Range("M65").Value = -14513.3335
$ws.Range("N65").Value = -30533
$ws.Range("H76").Value = 3816.5
$ws.Range("J76").Value = 3333
$ws.Range("L76").Value = 3333
$ws.Range("N76").Value = -3963
$ws.Range("H79").Value = 3816.5
$ws.Range("J79").Value = 3333
$ws.Range("L79").Value = 3333
$ws.Range("N79").Value = -5517
$ws.Range("H107").Value = 3477.818
$ws.Range("I107").Value = 2198.625
$ws.Range("J107").Value = 6889
$ws.Range("K107").Value = 2198.625
$ws.Range("L107").Value = 6889
$ws.Range("M107").Value = -278.625
$ws.Range("N107").Value = -10729
$ws.Range("H133").Value = 85710
$ws.Range("J133").Value = 85710
$ws.Range("L133").Value = 85710
$ws.Range("N133").Value = -95830
$ws.Range("H137").Value = 5597.9414
$ws.Range("J137").Value = 6324.636
$ws.Range("L137").Value = 18973.908
$ws.Range("N137").Value = -24073.908
$ws.Range("H138").Value = 7719.4243
$ws.Range("J138").Value = 8569.875
$ws.Range("L138").Value = 25709.625
$ws.Range("N138").Value = -35989.625
$ws.Range("H140").Value = 60516.668
$ws.Range("J140").Value = 59331.25
$ws.Range("L140").Value = 59331.25
$ws.Range("N140").Value = -69691.25
$ws.Range("H141").Value = 7089.2
$ws.Range("I141").Value = 7984.5713
$ws.Range("K141").Value = 23953.7139
$ws.Range("M141").Value = -18773.7139

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 832949
$ws.Range("I2").Value = 1343437.9
$ws.Range("J2").Value = 3404.5
$ws.Range("K2").Value = 1343437.9
$ws.Range("L2").Value = 3404.5
$ws.Range("M2").Value = -1343324.9
$ws.Range("N2").Value = -3630.5
$ws.Range("H63").Value = 2842.1428
$ws.Range("I63").Value = 2842.1428
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 2842.1428
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -2156.1428
$ws.Range("N63").ClearContents()
$ws.Range("H66").Value = 2842.1428
$ws.Range("I66").Value = 2842.1428
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 14210.714
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -10778.714
$ws.Range("N66").ClearContents()
$ws.Range("H116").Value = 832949
$ws.Range("I116").Value = 1343437.9
$ws.Range("J116").Value = 3404.5
$ws.Range("K116").Value = 1343437.9
$ws.Range("L116").Value = 3404.5
$ws.Range("M116").Value = -1341143.9
$ws.Range("N116").Value = -7992.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 832949
$ws.Range("I3").Value = 1343437.9
$ws.Range("J3").Value = 3404.5
$ws.Range("K3").Value = 1343437.9
$ws.Range("L3").Value = 3404.5
$ws.Range("M3").Value = -1343323.9
$ws.Range("N3").Value = -3632.5
$ws.Range("H61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("N61").ClearContents()
$ws.Range("H99").Value = 582205
$ws.Range("I99").Value = 1604476
$ws.Range("K99").Value = 1604476
$ws.Range("M99").Value = -1602978

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H55").Value = 5000
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 5000
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 5000
$ws.Range("M55").ClearContents()
$ws.Range("N55").Value = -5630
$ws.Range("H99").Value = 19656.666
$ws.Range("I99").Value = 32227.5
$ws.Range("K99").Value = 32227.5
$ws.Range("M99").Value = -30729.5
$ws.Range("H105").Value = 2842183.2
$ws.Range("I105").Value = 7576555.5
$ws.Range("K105").Value = 7576555.5
$ws.Range("M105").Value = -7574808.5
$ws.Range("H126").Value = 19656.666
$ws.Range("I126").Value = 32227.5
$ws.Range("K126").Value = 96682.5
$ws.Range("M126").Value = -94212.5
$ws.Range("H141").Value = 373695
$ws.Range("J141").Value = 373695
$ws.Range("L141").Value = 373695
$ws.Range("N141").Value = -384055

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 841.7727
$ws.Range("I122").Value = 689.8570999999999
$ws.Range("K122").Value = 6208.7139
$ws.Range("M122").Value = -3758.7139

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 530154.5
$ws.Range("I122").Value = 738882.9399999999
$ws.Range("K122").Value = 2216648.82
$ws.Range("M122").Value = -2214198.82

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 8337119.5
$ws.Range("I40").Value = 12502431
$ws.Range("K40").Value = 12502431
$ws.Range("M40").Value = -12502295
$ws.Range("H45").Value = 30000
$ws.Range("I45").Value = 30000
$ws.Range("K45").Value = 30000
$ws.Range("M45").Value = -29593
$ws.Range("H138").Value = 49998.332
$ws.Range("J138").Value = 49998.332
$ws.Range("L138").Value = 49998.332
$ws.Range("N138").Value = -60278.332
$ws.Range("H140").Value = 64692.6
$ws.Range("J140").Value = 64692.6
$ws.Range("L140").Value = 64692.6
$ws.Range("N140").Value = -75052.60000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 22735638
$ws.Range("I132").Value = 9336.111000000001
$ws.Range("K132").Value = 28008.333
$ws.Range("M132").Value = -25478.333
$ws.Range("H136").Value = 8711.706
$ws.Range("I136").Value = 1049.75
$ws.Range("J136").Value = 9733.299999999999
$ws.Range("K136").Value = 3149.25
$ws.Range("L136").Value = 29199.9
$ws.Range("M136").Value = -599.25
$ws.Range("N136").Value = -34299.89999999999
$ws.Range("H141").Value = 84000
$ws.Range("J141").Value = 84000
$ws.Range("L141").Value = 84000
$ws.Range("N141").Value = -94360
